# Capital costs / TODO sheet fix:
# - Mark "PV self-shading testing" (row 34) as Done.
# - Insert a new "Not done" TODO item right after it (new row 35):
#     What: "Re-arrange self-shading inputs in UI with system design? Check
#            inputs for usability in SDK"
#     Who:  Janine
#   This pushes all the following rows down by one (old row 35 -> 36, ...,
#   old row 61 -> 62), which Excel does automatically together with the
#   dependent formula (H17) and the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 36 ("Implementation of IEC
# 61853 algorithms in C++"), i.e. at position 35, right after the
# "PV self-shading testing" row. Excel shifts rows 35:61 down to 36:62,
# carries formatting from the row above, and auto-adjusts the SUM(D17:D39)
# formula in H17 to SUM(D17:D40) plus the sheet dimension to A1:I62.
$ws.Rows.Item(35).Insert()

# Row 34 ("PV self-shading testing") is now Done.
$ws.Range("A34").Value = "Done"

# Fill in the newly inserted row 35.
$ws.Range("A35").Value = "Not done"
$ws.Range("B35").Value = "Re-arrange self-shading inputs in UI with system design? Check inputs for usability in SDK"
$ws.Range("C35").Value = "Janine"

# Restore the view's active cell/selection to B36 (where the new row
# ended up moving the cursor in the original edit).
$ws.Range("B36").Select()
